$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 922.1667
$ws.Range("I111").Value = 914.2143
$ws.Range("J111").Value = 950
$ws.Range("K111").Value = 2742.6429
$ws.Range("L111").Value = 2850
$ws.Range("M111").Value = 324.3571000000002
$ws.Range("N111").Value = -8984
$ws.Range("H116").Value = 3441.25
$ws.Range("I116").Value = 2785.625
$ws.Range("J116").Value = 4315.4165
$ws.Range("K116").Value = 2785.625
$ws.Range("L116").Value = 4315.4165
$ws.Range("M116").Value = 656.375
$ws.Range("N116").Value = -11199.4165
$ws.Range("H132").Value = 3227525.5
$ws.Range("I132").Value = 3510146.2
$ws.Range("J132").Value = 5649.8
$ws.Range("K132").Value = 10530438.6
$ws.Range("L132").Value = 16949.4
$ws.Range("M132").Value = -10527908.6
$ws.Range("N132").Value = -22009.4
$ws.Range("H137").Value = 3225.2195
$ws.Range("I137").Value = 3844.524
$ws.Range("J137").Value = 2574.95
$ws.Range("K137").Value = 11533.572
$ws.Range("L137").Value = 7724.849999999999
$ws.Range("M137").Value = -8983.572
$ws.Range("N137").Value = -12824.85
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14707840
$ws.Range("I2").Value = 50001356
$ws.Range("K2").Value = 50001356
$ws.Range("M2").Value = -50001243
$ws.Range("H23").Value = 27125.75
$ws.Range("J23").Value = 21000
$ws.Range("L23").Value = 21000
$ws.Range("N23").Value = -21518
$ws.Range("H32").Value = 4334.9414
$ws.Range("I32").Value = 3436.5088
$ws.Range("J32").Value = 8990.454
$ws.Range("K32").Value = 3436.5088
$ws.Range("L32").Value = 8990.454
$ws.Range("M32").Value = -3149.5088
$ws.Range("N32").Value = -9564.454
$ws.Range("H102").Value = 2634.8276
$ws.Range("I102").Value = 2323.4614
$ws.Range("K102").Value = 2323.4614
$ws.Range("M102").Value = -701.4614000000001
$ws.Range("H116").Value = 14707840
$ws.Range("I116").Value = 50001356
$ws.Range("K116").Value = 50001356
$ws.Range("M116").Value = -49999062
$ws.Range("H122").Value = 3391.6191
$ws.Range("I122").Value = 2414.9333
$ws.Range("J122").Value = 5833.3335
$ws.Range("K122").Value = 7244.7999
$ws.Range("L122").Value = 17500.0005
$ws.Range("M122").Value = -4794.7999
$ws.Range("N122").Value = -22400.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14707840
$ws.Range("I3").Value = 50001356
$ws.Range("K3").Value = 50001356
$ws.Range("M3").Value = -50001242
$ws.Range("H107").Value = 1620.8667
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 2368.111
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 2368.111
$ws.Range("M107").Value = 1420
$ws.Range("N107").Value = -6208.111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2653.2
$ws.Range("I16").Value = 1974.75
$ws.Range("J16").Value = 2899.9092
$ws.Range("K16").Value = 1974.75
$ws.Range("L16").Value = 2899.9092
$ws.Range("M16").Value = -1687.75
$ws.Range("N16").Value = -3473.9092
$ws.Range("H31").Value = 2228.574
$ws.Range("I31").Value = 1339.6285
$ws.Range("K31").Value = 1339.6285
$ws.Range("M31").Value = -1044.6285
$ws.Range("H34").Value = 2228.574
$ws.Range("I34").Value = 1339.6285
$ws.Range("K34").Value = 1339.6285
$ws.Range("M34").Value = -1137.6285
$ws.Range("H50").Value = 11500
$ws.Range("I50").Value = 11500
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 11500
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -10875
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472
$ws.Range("H59").Value = 25975
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 27542.857
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 27542.857
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -29832.857
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696
$ws.Range("H99").Value = 2188.6
$ws.Range("I99").Value = 1293.5
$ws.Range("J99").Value = 2514.0908
$ws.Range("K99").Value = 1293.5
$ws.Range("L99").Value = 2514.0908
$ws.Range("M99").Value = 204.5
$ws.Range("N99").Value = -5510.0908
$ws.Range("H113").Value = 2653.2
$ws.Range("I113").Value = 1974.75
$ws.Range("J113").Value = 2899.9092
$ws.Range("K113").Value = 1974.75
$ws.Range("L113").Value = 2899.9092
$ws.Range("M113").Value = 195.25
$ws.Range("N113").Value = -7239.9092
$ws.Range("H126").Value = 2188.6
$ws.Range("I126").Value = 1293.5
$ws.Range("J126").Value = 2514.0908
$ws.Range("K126").Value = 3880.5
$ws.Range("L126").Value = 7542.2724
$ws.Range("M126").Value = -1410.5
$ws.Range("N126").Value = -12482.2724
$ws.Range("H134").Value = 1963.1936
$ws.Range("I134").Value = 1217.037
$ws.Range("K134").Value = 3651.111
$ws.Range("M134").Value = -1116.111
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1780
$ws.Range("J86").Value = 1780
$ws.Range("L86").Value = 5340
$ws.Range("N86").Value = -7712
$ws.Range("H89").Value = 1780
$ws.Range("J89").Value = 1780
$ws.Range("L89").Value = 16020
$ws.Range("N89").Value = -27876
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 568.7368
$ws.Range("I2").Value = 754
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 754
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = -641
$ws.Range("N2").Value = -276
$ws.Range("H7").Value = 6713335.5
$ws.Range("I7").Value = 20000000
$ws.Range("J7").Value = 70003
$ws.Range("K7").Value = 20000000
$ws.Range("L7").Value = 70003
$ws.Range("M7").Value = -19999888
$ws.Range("N7").Value = -70227
$ws.Range("H8").Value = 6713335.5
$ws.Range("I8").Value = 20000000
$ws.Range("J8").Value = 70003
$ws.Range("K8").Value = 20000000
$ws.Range("L8").Value = 70003
$ws.Range("M8").Value = -19999861
$ws.Range("N8").Value = -70281
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1939.1
$ws.Range("I16").Value = 1095.75
$ws.Range("J16").Value = 5312.5
$ws.Range("K16").Value = 1095.75
$ws.Range("L16").Value = 5312.5
$ws.Range("M16").Value = -925.75
$ws.Range("N16").Value = -5652.5
$ws.Range("H46").Value = 1911
$ws.Range("I46").Value = 946.6667
$ws.Range("K46").Value = 946.6667
$ws.Range("M46").Value = -758.6667
$ws.Range("H122").Value = 2956
$ws.Range("I122").Value = 2470.75
$ws.Range("J122").Value = 3926.5
$ws.Range("K122").Value = 7412.25
$ws.Range("L122").Value = 11779.5
$ws.Range("M122").Value = -4962.25
$ws.Range("N122").Value = -16679.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3562.6365
$ws.Range("I122").Value = 2766.6667
$ws.Range("J122").Value = 4517.8
$ws.Range("K122").Value = 8300.000100000001
$ws.Range("L122").Value = 13553.4
$ws.Range("M122").Value = -5850.000100000001
$ws.Range("N122").Value = -18453.4
$ws.Range("H135").Value = 89269.164
$ws.Range("J135").Value = 89269.164
$ws.Range("L135").Value = 89269.164
$ws.Range("N135").Value = -99409.164
$ws.Range("H136").Value = 2489
$ws.Range("I136").Value = 982
$ws.Range("K136").Value = 2946
$ws.Range("M136").Value = -396
$ws.Range("H141").Value = 26974.092
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 26974.092
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 26974.092
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -37334.092
